$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.862146333333333
$ws.Range("H2").Value = 8.586439
$ws.Range("I2").Value = 0.9778268096017091
$ws.Range("J2").Value = 0.9778268096017091
$ws.Range("Q2").Value = 0.1280343000265556
$ws.Range("R2").Value = 1.152308700239
$ws.Range("S2").Value = 0.9778268096017091
$ws.Range("T2").Value = 0.9778268096017091

# Row 3 updates
$ws.Range("I3").Value = 0.02217319039829088
$ws.Range("J3").Value = 0.02217319039829088
$ws.Range("S3").Value = 0.02217319039829088
$ws.Range("T3").Value = 0.02217319039829088
